$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 - DKS / Dekstop / Desktop Computer / eng
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"
$ws.Range("D2").Value = "eng"

# Update row 3 - DKS / Arabic / Arabic / ara
$ws.Range("A3").Value = "DKS"
$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"
$ws.Range("D3").Value = "ara"

# Update row 4 - DKS / Ordinateur / Ordinateurs de bureau / fra
$ws.Range("A4").Value = "DKS"
$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"
$ws.Range("D4").Value = "fra"

# Set the active cell selection to D10
$ws.Range("D10").Select()

# Set up page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
